# Row 7 and Row 8 had their record data swapped (the two species
# observations traded places in the sheet), except column AW
# ("Rapportör") which remained "David Isaksson" for both rows.
#
# Capture the original values first (both rows), then write the swapped
# values back so the script is resilient to ordering.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Original row 7 values --------------------------------------------
$A7  = 131106436
$B7  = 5493
$E7  = 101410
$F7  = "Reliktbock"
$G7  = "Nothorhina muricata"
$H7  = "(Dalman, 1817)"
$I7  = "2"
$J7  = "ex."
$P7  = "Svartmyran, Mpd"
$Q7  = 616762
$R7  = 6934714
$X7  = "2025_0743"
$Z7  = "11:39"
$AB7 = "11:39"
$AC7 = "Två kläckhål"
$AX7 = "David Isaksson"

# --- Original row 8 values --------------------------------------------
$A8  = 131108352
$B8  = 80214
$E8  = 388
$F8  = "Stiftgelélav"
$G8  = "Collema furfuraceum"
$H8  = "(Arnold) Du Rietz"
$I8  = "1"
$J8  = "bålar"
$P8  = "S Svartmyran, Mpd"
$Q8  = 616863
$R8  = 6934788
$X8  = "2025_0758"
$Z8  = "14:47"
$AB8 = "14:47"
$AC8 = "Asp"
$AX8 = "Måns Svensson"

# --- Write row 7 with what used to be row 8's data ---------------------
$ws.Range("A7").Value = $A8
$ws.Range("B7").Value = $B8
$ws.Range("E7").Value = $E8
$ws.Range("F7").Value = $F8
$ws.Range("G7").Value = $G8
$ws.Range("H7").Value = $H8
$ws.Range("I7").Value = "'" + $I8
$ws.Range("J7").Value = $J8
$ws.Range("P7").Value = $P8
$ws.Range("Q7").Value = $Q8
$ws.Range("R7").Value = $R8
$ws.Range("X7").Value = $X8
$ws.Range("Z7").Value = $Z8
$ws.Range("AB7").Value = $AB8
$ws.Range("AC7").Value = $AC8
$ws.Range("AX7").Value = $AX8

# --- Write row 8 with what used to be row 7's data ---------------------
$ws.Range("A8").Value = $A7
$ws.Range("B8").Value = $B7
$ws.Range("E8").Value = $E7
$ws.Range("F8").Value = $F7
$ws.Range("G8").Value = $G7
$ws.Range("H8").Value = $H7
$ws.Range("I8").Value = "'" + $I7
$ws.Range("J8").Value = $J7
$ws.Range("P8").Value = $P7
$ws.Range("Q8").Value = $Q7
$ws.Range("R8").Value = $R7
$ws.Range("X8").Value = $X7
$ws.Range("Z8").Value = $Z7
$ws.Range("AB8").Value = $AB7
$ws.Range("AC8").Value = $AC7
$ws.Range("AX8").Value = $AX7
